$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "56.400.91"
$ws.Range("E2").Value = "  +3.38%  "
$ws.Range("D3").Value = "2.497.40"
$ws.Range("E3").Value = "  +2.48%  "
$ws.Range("E4").Value = "  +0.26%  "
$ws.Range("D5").Value = "'489.13"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  +4.38%  "
$ws.Range("D6").Value = "'146.57"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = "  +10.60%  "
$ws.Range("D7").Value = "'0.996"
$ws.Range("D7").ClearFormats()
$ws.Range("E7").Value = "  +0.08%  "
$ws.Range("E8").Value = "  +4.95%  "
$ws.Range("D9").Value = "2.517.57"
$ws.Range("E9").Value = "  +2.62%  "
$ws.Range("E10").Value = "  +8.31%  "
$ws.Range("D11").Value = "'0.0977"
$ws.Range("D11").ClearFormats()
$ws.Range("D12").Value = "'0.333"
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = "  +4.31%  "
$ws.Range("E13").Value = "  +1.28%  "
$ws.Range("D14").Value = "2.943.36"
$ws.Range("E14").Value = "  +2.85%  "
$ws.Range("D15").Value = "56.459.81"
$ws.Range("E15").Value = "  +4.03%  "
$ws.Range("D16").Value = "'21.28"
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = "  +6.95%  "
$ws.Range("E17").Value = "  +2.37%  "
$ws.Range("D18").Value = "2.512.37"
$ws.Range("E18").Value = "  +2.15%  "
$ws.Range("D19").Value = "'4.53"
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = "  +7.25%  "
$ws.Range("D20").Value = "'10.23"
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = "  +8.33%  "
$ws.Range("D21").Value = "'321.01"
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = "  +2.00%  "
$ws.Range("D22").Value = "'0.998"
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = "  +0.86%  "
$ws.Range("E23").Value = "  +8.58%  "
$ws.Range("D24").Value = "'58.72"
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = "  +3.32%  "
$ws.Range("E25").Value = "  +6.66%  "
$ws.Range("D26").Value = "'0.166"
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = "  +8.38%  "
$ws.Range("D27").Value = "'0.997"
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = "  -0.93%  "
$ws.Range("D28").Value = "2.618.94"
$ws.Range("E28").Value = "  +3.02%  "
$ws.Range("D29").Value = "'7.63"
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = "  +5.28%  "
$ws.Range("D30").Value = "0.0₃0794"
$ws.Range("E30").Value = "  +8.90%  "
$ws.Range("E31").Value = "  +0.30%  "
$ws.Range("D32").Value = "'148.70"
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = "  -0.94%  "
$ws.Range("D33").Value = "'18.29"
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = "  +2.54%  "
$ws.Range("E34").Value = "  +4.81%  "
$ws.Range("D35").Value = "'5.22"
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = "  +3.46%  "
$ws.Range("E36").Value = "  +7.57%  "
$ws.Range("E37").Value = "  +4.99%  "
$ws.Range("D38").Value = "'0.869"
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = "  +7.82%  "
$ws.Range("D39").Value = "'34.23"
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = "  +1.58%  "
$ws.Range("D40").Value = "'3.54"
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = "  +7.60%  "
$ws.Range("E41").Value = "  +2.89%  "
$ws.Range("E42").Value = "  +5.29%  "
$ws.Range("D43").Value = "'0.994"
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = "  -0.05%  "
$ws.Range("E44").Value = "  +7.52%  "
$ws.Range("D45").Value = "'4.87"
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = "  +12.56%  "
$ws.Range("D46").Value = "'261.38"
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = "  +16.69%  "
$ws.Range("E47").Value = "  +3.79%  "
$ws.Range("D48").Value = "'10.19"
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = "  -0.23%  "
$ws.Range("D49").Value = "'0.0911"
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = "  +4.11%  "
$ws.Range("D50").Value = "1.917.23"
$ws.Range("E50").Value = "  -2.01%  "
$ws.Range("D51").Value = "'17.71"
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = "  +5.74%  "
